# edit.ps1 - Applies the Aug 11 2023 04:30:48 UTC GitHub Actions crypto data refresh
# to the cryptos worksheet (Coin, Link, Price, Volume(1h) columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.375.33'
$ws.Range('E2').Value = '  -0.50%  '

$ws.Range('D3').Value = '1.844.67'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('D4').Value = '''0.9982'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').Value = '''240.45'
$ws.Range('E5').Value = '  -1.07%  '

$ws.Range('D6').Value = '''0.6390'
$ws.Range('E6').Value = '  +0.38%  '

$ws.Range('D7').Value = '''0.9996'
$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '''0.07556'
$ws.Range('E8').Value = '  -0.23%  '

$ws.Range('D9').Value = '''0.2967'
$ws.Range('E9').Value = '  -1.10%  '

$ws.Range('D10').Value = '''24.76'
$ws.Range('E10').Value = '  +2.15%  '

$ws.Range('E11').Value = '  +0.52%  '

$ws.Range('D12').Value = '1.867.69'
$ws.Range('E12').Value = '  +0.47%  '

$ws.Range('D13').Value = '''4.991'
$ws.Range('E13').Value = '  -0.90%  '

$ws.Range('D14').Value = '''0.6845'
$ws.Range('E14').Value = '  -0.47%  '

$ws.Range('D15').Value = '''83.21'
$ws.Range('E15').Value = '  -0.95%  '

$ws.Range('D16').Value = '''0.000009956'
$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('D17').Value = '''6.176'
$ws.Range('E17').Value = '  -1.52%  '

$ws.Range('D18').Value = '29.396.38'
$ws.Range('E18').Value = '  -0.56%  '

$ws.Range('D19').Value = '''229.70'
$ws.Range('E19').Value = '  -3.17%  '

$ws.Range('E20').Value = '  -0.63%  '

$ws.Range('D21').Value = '''0.9995'
$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('D22').Value = '''7.565'
$ws.Range('E22').Value = '  -0.76%  '

$ws.Range('D23').Value = '''0.9997'
$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '''155.99'
$ws.Range('E24').Value = '  -0.42%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '''0.1406'
$ws.Range('E25').Value = '  +0.79%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '''8.393'
$ws.Range('E26').Value = '  -0.76%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''17.67'
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = '''1.465'
$ws.Range('E28').Value = '  -1.65%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '''0.05720'
$ws.Range('E29').Value = '  -2.91%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''1.247'
$ws.Range('E30').Value = '  -2.39%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''4.134'
$ws.Range('E31').Value = '  +0.10%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''4.034'
$ws.Range('E32').Value = '  -0.87%  '

$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = '''1.848'
$ws.Range('E33').Value = '  -2.71%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''1.157'
$ws.Range('E34').Value = '  -1.47%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.7171'
$ws.Range('E35').Value = '  -0.43%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '''2.589'
$ws.Range('E36').Value = '  -0.40%  '

$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.254.54'
$ws.Range('E37').Value = '  +2.10%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01811'
$ws.Range('E38').Value = '  +1.78%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.787'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''0.9106'
$ws.Range('E40').Value = '  -0.27%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''0.9994'
$ws.Range('E41').Value = '  -0.03%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '''101.60'
$ws.Range('E42').Value = '  -0.35%  '

$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''66.48'
$ws.Range('E43').Value = '  -1.46%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = '''7.064'
$ws.Range('E44').Value = '  -4.49%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '''0.00000000118'
$ws.Range('E45').Value = '  +0.31%  '

$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Value = '''0.4030'
$ws.Range('E46').Value = '  -0.40%  '

$ws.Range('D47').Value = '''9.132'
$ws.Range('E47').Value = '  -0.13%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.708'
$ws.Range('E48').Value = '  +0.62%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.1128'
$ws.Range('E49').Value = '  +0.33%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05738'
$ws.Range('E50').Value = '  -0.18%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.4627'
$ws.Range('E51').Value = '  -0.21%  '
